# Add a new "ano" (year) column in column A, with header in A2 and
# value 2023 for every data row (A3:A35), matching the B/C columns that
# already span rows 2-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("A2").Value = "ano"

# Fill A3:A35 with the year value 2023
$ws.Range("A3:A35").Value = 2023

# Update the selection to match the edited range (A4:A35, active cell A4)
$ws.Range("A4:A35").Select()
